$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of quotes for 2025-09-12 (serial date 45912)
$ws.Cells.Item(8, 1).Value = 45912
$ws.Cells.Item(8, 1).NumberFormat = $ws.Cells.Item(7, 1).NumberFormat

$ws.Cells.Item(8, 2).Value = "21,1936"
$ws.Cells.Item(8, 3).Value = "14,9727"
$ws.Cells.Item(8, 4).Value = "14,8984"
$ws.Cells.Item(8, 5).Value = "14,8984"
